$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Replace "Good Morning" with "GIT UPDATE" in cell E8, and select that cell
# (mirrors the author editing the cell directly in Excel and leaving it selected)
$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
